$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new reservation (Gerlinde Weiss) was inserted as row 48, pushing the
# existing rows 48-50 (Franziska Lindermeier, Claudine Fleury, TOTAL) down
# to rows 49-51.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new reservation data
$ws.Cells.Item(48, 1).Value = "Gerlinde Weiss"
$ws.Cells.Item(48, 2).Value = "Booking"

# Phone number must stay text (keeps the leading "+")
$ws.Cells.Item(48, 3).NumberFormat = "@"
$ws.Cells.Item(48, 3).Value = "+4369912047111"
$ws.Cells.Item(48, 3).ClearFormats()

# date_arrivee / date_depart (2026-01-02 -> 2026-01-03), same date style
# used by the other rows in this column
$ws.Cells.Item(48, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(48, 4).Value = 46024
$ws.Cells.Item(48, 5).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(48, 5).Value = 46025

$ws.Cells.Item(48, 6).Value = 1
$ws.Cells.Item(48, 7).Value = 535.24
$ws.Cells.Item(48, 8).Value = 439.86
$ws.Cells.Item(48, 9).Value = 95.38
$ws.Cells.Item(48, 10).Value = 17.82
$ws.Cells.Item(48, 11).Value = 2026
$ws.Cells.Item(48, 12).Value = 1

$ws.Cells.Item(48, 13).ClearFormats()
$ws.Cells.Item(48, 14).ClearFormats()
$ws.Cells.Item(48, 15).ClearFormats()
